$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix B5: was stored as text, now stored as a genuine number
$ws.Range("B5").Value = 23087278

# Row 6 - new record
$ws.Range("A6").Value = "2025-10-28 13:24:32"
$ws.Range("B6").Value = 23153422
$ws.Range("C6").Value = "Daniel "
$ws.Range("D6").Value = "Ejecutado en Campo"
$ws.Range("E6").Value = "23153422_1_20251028_132432.pdf"
$ws.Range("F6").Value = "Sin imágenes"
$ws.Range("G6").Value = "LIZANA PATRICIA BEDOYA MEJIA"
$ws.Range("H6").Value = "CR 1 ESTE CL 47 C -18"
$ws.Range("I6").Value = "VENCIDO"
$ws.Range("J6").Value = 1035861667
$ws.Range("K6").Value = "Formulario"

# Row 7 - new record (pedido stays text here)
$ws.Range("A7").Value = "2025-10-28 13:39:07"
$ws.Range("B7").NumberFormat = "@"
$ws.Range("B7").Value = "23308435"
$ws.Range("B7").Style = "Normal"
$ws.Range("C7").Value = "descarte 440"
$ws.Range("D7").Value = "Descartado"
$ws.Range("E7").Value = "23308435_1_20251028_133907.pdf"
$ws.Range("F7").Value = "Sin imágenes"
$ws.Range("G7").Value = "JESUS AMADO JIMENEZ SUCERQUIA"
$ws.Range("H7").Value = "CR 21 CL 56 BD -5 (INTERIOR 201 )"
$ws.Range("I7").Value = "VENCIDO"
$ws.Range("J7").Value = 71263185
$ws.Range("K7").Value = "Formulario"
